$wb = $excel.ActiveWorkbook

# Sheet "展览" - update F column ("想去人数") values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 61
$wsExhibit.Range("F6").Value = 35
$wsExhibit.Range("F8").Value = 90
$wsExhibit.Range("F9").Value = 8581
$wsExhibit.Range("F12").Value = 1138
$wsExhibit.Range("F13").Value = 935
$wsExhibit.Range("F14").Value = 93
$wsExhibit.Range("F17").Value = 216
$wsExhibit.Range("F19").Value = 229
$wsExhibit.Range("F20").Value = 978

# Sheet "全部类型" - update F column ("想去人数") values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 61
$wsAll.Range("F7").Value = 35
$wsAll.Range("F10").Value = 90
$wsAll.Range("F11").Value = 8581
$wsAll.Range("F14").Value = 1138
$wsAll.Range("F15").Value = 935
$wsAll.Range("F16").Value = 93
$wsAll.Range("F19").Value = 216
$wsAll.Range("F21").Value = 229
$wsAll.Range("F22").Value = 978
